$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new quotation row (row 18) for date 2025-09-22 (Excel serial 45922)
$ws.Range("A18").Value = 45922

# Copy formatting (number format / style) from the row above so the new
# date cell renders the same way as the rest of column A
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B18").Value = "14,9404"
$ws.Range("C18").Value = "15,0881"
$ws.Range("D18").Value = "14,9404"
$ws.Range("E18").Value = "14,9404"
